# Scheduled-runner price refresh: update cached market/profit figures on
# several leve rows across the Pandaemonium crafting-job sheets.
# (H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)

$wb = $excel.ActiveWorkbook

function Set-RowValues($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $v = $Values[$col]
        $cell = $ws.Range("$col$Row")
        if ($null -eq $v) {
            $cell.Value = ""
        } else {
            $cell.Value = $v
        }
    }
}

# --- ALC ---
Set-RowValues "ALC" 28 @{
    H = 579.37036; I = 356.58823; J = 958.1; K = 356.58823; L = 958.1; M = 128.41177; N = -1928.1
}
Set-RowValues "ALC" 62 @{
    H = 1816.8182; I = 1746.6666; J = 1901; K = 1746.6666; L = 1901; M = -1122.6666; N = -3149
}
Set-RowValues "ALC" 64 @{
    H = 3320.6216; I = 3166.5715; J = 3799.889; K = 3166.5715; L = 3799.889; M = -2918.5715; N = -4295.889
}
Set-RowValues "ALC" 65 @{
    H = 1816.8182; I = 1746.6666; J = 1901; K = 8733.333000000001; L = 9505; M = -5613.333000000001; N = -15745
}
Set-RowValues "ALC" 67 @{
    H = 3320.6216; I = 3166.5715; J = 3799.889; K = 3166.5715; L = 3799.889; M = -2308.5715; N = -5515.889
}
Set-RowValues "ALC" 111 @{
    H = 4443.5; I = 6927; J = 1960; K = 20781; L = 5880; M = -17714; N = -12014
}
Set-RowValues "ALC" 115 @{
    H = 2887.111; I = 2880.6667; J = 2900; K = 8642.000100000001; L = 8700; M = -7075.000100000001; N = -11834
}
Set-RowValues "ALC" 133 @{
    H = 50063.35; J = 50063.35; L = 50063.35; N = -60183.35
}
Set-RowValues "ALC" 134 @{
    H = 52828.75; J = 52828.75; L = 52828.75; N = -62968.75
}
Set-RowValues "ALC" 137 @{
    H = 648650.75; I = 1731.48; J = 1547149.8; K = 5194.440000000001; L = 4641449.4; M = -2644.440000000001; N = -4646549.4
}

# --- ARM ---
Set-RowValues "ARM" 2 @{
    H = 1173.9375; I = 1201; K = 1201; M = -1088
}
Set-RowValues "ARM" 4 @{
    H = 365.33334; I = 365.33334; K = 365.33334; M = -249.33334
}
Set-RowValues "ARM" 110 @{
    H = 1121.8334; I = 882.75; J = 1600; K = 882.75; L = 1600; M = 1162.25; N = -5690
}
Set-RowValues "ARM" 116 @{
    H = 1173.9375; I = 1201; K = 1201; M = 1093
}
Set-RowValues "ARM" 122 @{
    H = 5684771.5; I = 3276.4666; K = 9829.399800000001; M = -7379.399800000001
}

# --- BSM ---
Set-RowValues "BSM" 3 @{
    H = 1173.9375; I = 1201; K = 1201; M = -1087
}
Set-RowValues "BSM" 29 @{
    H = 1270.6666; I = 1270.6666; K = 1270.6666; M = -981.6666
}
Set-RowValues "BSM" 36 @{
    H = 10356; I = 3134; J = 24800; K = 3134; L = 24800; M = -2600; N = -25868
}

# --- CRP ---
Set-RowValues "CRP" 7 @{
    H = 43.333332; I = 43.333332; J = 0; K = 43.333332; L = 0; M = 69.666668; N = $null
}
Set-RowValues "CRP" 99 @{
    H = 1632.75; I = 1692.4; J = 1533.3334; K = 1692.4; L = 1533.3334; M = -194.4000000000001; N = -4529.3334
}
Set-RowValues "CRP" 107 @{
    H = 404.83334; I = 401.57895; J = 417.2; K = 401.57895; L = 417.2; M = 1518.42105; N = -4257.2
}
Set-RowValues "CRP" 126 @{
    H = 1632.75; I = 1692.4; J = 1533.3334; K = 5077.200000000001; L = 4600.0002; M = -2607.200000000001; N = -9540.0002
}
Set-RowValues "CRP" 132 @{
    H = 2325.6296; I = 1874.5714; J = 3904.3333; K = 5623.7142; L = 11712.9999; M = -3093.7142; N = -16772.9999
}

# --- CUL ---
Set-RowValues "CUL" 41 @{
    H = 320.4; I = 150.5; J = 1000; K = 451.5; L = 3000; M = -113.5; N = -3676
}

# --- GSM ---
Set-RowValues "GSM" 133 @{
    H = 0; J = 0; L = 0; N = $null
}

# --- WVR ---
Set-RowValues "WVR" 107 @{
    H = 1131.5416; I = 366.5; J = 1514.0625; K = 1099.5; L = 4542.1875; M = 820.5; N = -8382.1875
}
Set-RowValues "WVR" 136 @{
    H = 4836.6274; I = 2711.1365; J = 6449.069; K = 8133.4095; L = 19347.207; M = -5583.4095; N = -24447.207
}
